# Insert a new data row before the current row 172, shifting the existing
# rows 172-272 down to 173-273 (their values stay the same, only their row
# numbers change). Populate the newly inserted row 172 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 172 (pushes old row 172.. down by one row).
$ws.Rows.Item(172).Insert()

# Fill in the values for the new row 172.
$ws.Cells.Item(172, 1).Value = 7
$ws.Cells.Item(172, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(172, 3).Value = "Ñuble"
$ws.Cells.Item(172, 4).Value = 44438
$ws.Cells.Item(172, 5).Value = 16
$ws.Cells.Item(172, 6).Value = 100112020
$ws.Cells.Item(172, 7).Value = "Tomate"
$ws.Cells.Item(172, 8).Value = "Larga vida"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 300
$ws.Cells.Item(172, 11).Value = 5500
$ws.Cells.Item(172, 12).Value = 6000
$ws.Cells.Item(172, 13).Value = 5750
$ws.Cells.Item(172, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(172, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(172, 16).Value = 575
$ws.Cells.Item(172, 17).Value = 10
$ws.Cells.Item(172, 18).Value = "Hortaliza"
